$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new A:C data block (rows 2..113) as a 2D array for a single bulk write.
$data = New-Object 'object[,]' 112,3
$data[0,0] = 0
$data[0,1] = 0
$data[0,2] = 20.46315725008776
$data[1,0] = 1
$data[1,1] = 1.1
$data[1,2] = -8.871362001440087
$data[2,0] = 2
$data[2,1] = 2.2
$data[2,2] = 6.803275573276045
$data[3,0] = 3
$data[3,1] = 3.3
$data[3,2] = -2.753470682257613
$data[4,0] = 4
$data[4,1] = 4.4
$data[4,2] = 12.74719179615108
$data[5,0] = 5
$data[5,1] = 5.5
$data[5,2] = -8.996605711510755
$data[6,0] = 6
$data[6,1] = 6.6
$data[6,2] = -4.853370001023103
$data[7,0] = 7
$data[7,1] = 7.699999999999999
$data[7,2] = -14.01748220980055
$data[8,0] = 8
$data[8,1] = 8.799999999999999
$data[8,2] = 4.527168444763065
$data[9,0] = 9
$data[9,1] = 9.899999999999999
$data[9,2] = -5.391267436463464
$data[10,0] = 10
$data[10,1] = 11
$data[10,2] = 9.690016386984103
$data[11,0] = 11
$data[11,1] = 12.1
$data[11,2] = 5.939281940682445
$data[12,0] = 12
$data[12,1] = 13.2
$data[12,2] = -7.634711911220005
$data[13,0] = 13
$data[13,1] = 14.3
$data[13,2] = -6.263356606729044
$data[14,0] = 14
$data[14,1] = 15.4
$data[14,2] = 0.932446025522675
$data[15,0] = 15
$data[15,1] = 16.5
$data[15,2] = 0.7346028022933537
$data[16,0] = 16
$data[16,1] = 17.6
$data[16,2] = 14.63082831826653
$data[17,0] = 17
$data[17,1] = 18.7
$data[17,2] = 6.912665085636519
$data[18,0] = 18
$data[18,1] = 19.8
$data[18,2] = -2.946189457186951
$data[19,0] = 19
$data[19,1] = 20.9
$data[19,2] = -15.17831728192271
$data[20,0] = 20
$data[20,1] = 22
$data[20,2] = 9.34162203604075
$data[21,0] = 21
$data[21,1] = 23.1
$data[21,2] = 1.74485536124526
$data[22,0] = 22
$data[22,1] = 24.20000000000001
$data[22,2] = 9.679005220521654
$data[23,0] = 23
$data[23,1] = 25.30000000000001
$data[23,2] = -7.426001500536788
$data[24,0] = 24
$data[24,1] = 26.40000000000001
$data[24,2] = 10.98027960908465
$data[25,0] = 25
$data[25,1] = 27.50000000000001
$data[25,2] = -12.64345235274561
$data[26,0] = 26
$data[26,1] = 28.60000000000001
$data[26,2] = -10.17970735623242
$data[27,0] = 27
$data[27,1] = 29.70000000000001
$data[27,2] = -1.787385927220983
$data[28,0] = 28
$data[28,1] = 30.80000000000001
$data[28,2] = -14.82256042731548
$data[29,0] = 29
$data[29,1] = 31.90000000000002
$data[29,2] = 10.9447860140084
$data[30,0] = 30
$data[30,1] = 33.00000000000001
$data[30,2] = 3.656340843174286
$data[31,0] = 31
$data[31,1] = 34.10000000000002
$data[31,2] = -0.2563721846741766
$data[32,0] = 32
$data[32,1] = 35.20000000000002
$data[32,2] = -7.661025087283079
$data[33,0] = 33
$data[33,1] = 36.30000000000002
$data[33,2] = -18.73522161967288
$data[34,0] = 34
$data[34,1] = 37.40000000000002
$data[34,2] = -2.336404862798283
$data[35,0] = 35
$data[35,1] = 38.50000000000002
$data[35,2] = -6.117131279581152
$data[36,0] = 36
$data[36,1] = 39.60000000000002
$data[36,2] = -10.57595917320752
$data[37,0] = 37
$data[37,1] = 40.70000000000002
$data[37,2] = -5.504797896560647
$data[38,0] = 38
$data[38,1] = 41.80000000000003
$data[38,2] = -1.907095634737647
$data[39,0] = 39
$data[39,1] = 42.90000000000003
$data[39,2] = -4.054855616330562
$data[40,0] = 40
$data[40,1] = 44.00000000000003
$data[40,2] = 2.167677736431082
$data[41,0] = 41
$data[41,1] = 45.10000000000003
$data[41,2] = -14.61331937269333
$data[42,0] = 42
$data[42,1] = 46.20000000000003
$data[42,2] = 1.763213644053386
$data[43,0] = 43
$data[43,1] = 47.30000000000003
$data[43,2] = -7.409518265668426
$data[44,0] = 44
$data[44,1] = 48.40000000000003
$data[44,2] = 9.233460107240255
$data[45,0] = 45
$data[45,1] = 49.50000000000004
$data[45,2] = -3.10008851554772
$data[46,0] = 46
$data[46,1] = 50.60000000000004
$data[46,2] = 7.056147080660277
$data[47,0] = 47
$data[47,1] = 51.70000000000004
$data[47,2] = 4.311363025421969
$data[48,0] = 48
$data[48,1] = 52.80000000000004
$data[48,2] = -1.915697423997605
$data[49,0] = 49
$data[49,1] = 53.90000000000004
$data[49,2] = 7.278561227610139
$data[50,0] = 50
$data[50,1] = 55.00000000000004
$data[50,2] = 2.319948216946591
$data[51,0] = 51
$data[51,1] = 56.10000000000004
$data[51,2] = -2.786791244652302
$data[52,0] = 52
$data[52,1] = 57.20000000000005
$data[52,2] = -6.267557027072623
$data[53,0] = 53
$data[53,1] = 58.30000000000005
$data[53,2] = -19.75759761966443
$data[54,0] = 54
$data[54,1] = 59.40000000000005
$data[54,2] = 13.41625513447663
$data[55,0] = 55
$data[55,1] = 60.50000000000005
$data[55,2] = -7.91024938243501
$data[56,0] = 56
$data[56,1] = 61.60000000000005
$data[56,2] = -5.619680341410904
$data[57,0] = 57
$data[57,1] = 62.70000000000005
$data[57,2] = -5.631614546368065
$data[58,0] = 58
$data[58,1] = 63.80000000000005
$data[58,2] = -2.864770294520354
$data[59,0] = 59
$data[59,1] = 64.90000000000005
$data[59,2] = -2.165155831899343
$data[60,0] = 60
$data[60,1] = 66.00000000000004
$data[60,2] = 8.55312234478756
$data[61,0] = 61
$data[61,1] = 67.10000000000004
$data[61,2] = -10.50665060936482
$data[62,0] = 62
$data[62,1] = 68.20000000000003
$data[62,2] = -6.559906776314818
$data[63,0] = 63
$data[63,1] = 69.30000000000003
$data[63,2] = -13.73044460349132
$data[64,0] = 64
$data[64,1] = 70.40000000000002
$data[64,2] = -7.113486767825908
$data[65,0] = 65
$data[65,1] = 71.50000000000001
$data[65,2] = -0.6937881474173797
$data[66,0] = 66
$data[66,1] = 72.60000000000001
$data[66,2] = -4.783216995699267
$data[67,0] = 67
$data[67,1] = 73.7
$data[67,2] = 1.018119696560815
$data[68,0] = 68
$data[68,1] = 74.8
$data[68,2] = -5.532114454705477
$data[69,0] = 69
$data[69,1] = 75.89999999999999
$data[69,2] = -6.881708503923551
$data[70,0] = 70
$data[70,1] = 76.99999999999999
$data[70,2] = -10.06953039859076
$data[71,0] = 71
$data[71,1] = 78.09999999999998
$data[71,2] = 4.776997964709762
$data[72,0] = 72
$data[72,1] = 79.19999999999997
$data[72,2] = -10.317878486377
$data[73,0] = 73
$data[73,1] = 80.29999999999997
$data[73,2] = -10.17673111843692
$data[74,0] = 74
$data[74,1] = 81.39999999999996
$data[74,2] = 7.710533746143586
$data[75,0] = 75
$data[75,1] = 82.49999999999996
$data[75,2] = -1.507901070348139
$data[76,0] = 76
$data[76,1] = 83.59999999999995
$data[76,2] = -1.113437540682142
$data[77,0] = 77
$data[77,1] = 84.69999999999995
$data[77,2] = -4.208875242643325
$data[78,0] = 78
$data[78,1] = 85.79999999999994
$data[78,2] = -13.35408333101818
$data[79,0] = 79
$data[79,1] = 86.89999999999993
$data[79,2] = 2.438968470909062
$data[80,0] = 80
$data[80,1] = 87.99999999999993
$data[80,2] = -2.087409342912719
$data[81,0] = 81
$data[81,1] = 89.09999999999992
$data[81,2] = 6.747578769235247
$data[82,0] = 82
$data[82,1] = 90.19999999999992
$data[82,2] = -8.558197192471734
$data[83,0] = 83
$data[83,1] = 91.29999999999991
$data[83,2] = 5.806552822150973
$data[84,0] = 84
$data[84,1] = 92.39999999999991
$data[84,2] = -7.048718900101809
$data[85,0] = 85
$data[85,1] = 93.4999999999999
$data[85,2] = -16.09709705985068
$data[86,0] = 86
$data[86,1] = 94.59999999999989
$data[86,2] = -16.14303247086716
$data[87,0] = 87
$data[87,1] = 95.69999999999989
$data[87,2] = 1.449554180435693
$data[88,0] = 88
$data[88,1] = 96.79999999999988
$data[88,2] = 8.330247609231924
$data[89,0] = 89
$data[89,1] = 97.89999999999988
$data[89,2] = 8.58578527482352
$data[90,0] = 90
$data[90,1] = 98.99999999999987
$data[90,2] = 4.435036605073154
$data[91,0] = 91
$data[91,1] = 100.0999999999999
$data[91,2] = -6.914884225776671
$data[92,0] = 92
$data[92,1] = 101.1999999999999
$data[92,2] = 5.051340354079974
$data[93,0] = 93
$data[93,1] = 102.2999999999999
$data[93,2] = 12.88220840161638
$data[94,0] = 94
$data[94,1] = 103.3999999999998
$data[94,2] = -2.553188231157454
$data[95,0] = 95
$data[95,1] = 104.4999999999998
$data[95,2] = 11.89150616439616
$data[96,0] = 96
$data[96,1] = 105.5999999999998
$data[96,2] = 9.587071702997644
$data[97,0] = 97
$data[97,1] = 106.6999999999998
$data[97,2] = -5.414525401398792
$data[98,0] = 98
$data[98,1] = 107.7999999999998
$data[98,2] = -13.11977186947088
$data[99,0] = 99
$data[99,1] = 108.8999999999998
$data[99,2] = 0.7996134891376993
$data[100,0] = 100
$data[100,1] = 109.9999999999998
$data[100,2] = 9.818744208345191
$data[101,0] = 101
$data[101,1] = 111.0999999999998
$data[101,2] = 5.063682026918908
$data[102,0] = 102
$data[102,1] = 112.1999999999998
$data[102,2] = -11.75914295985498
$data[103,0] = 103
$data[103,1] = 113.2999999999998
$data[103,2] = -0.9397642029997761
$data[104,0] = 104
$data[104,1] = 114.3999999999998
$data[104,2] = -3.216744901127722
$data[105,0] = 105
$data[105,1] = 115.4999999999998
$data[105,2] = -11.20916337697438
$data[106,0] = 106
$data[106,1] = 116.5999999999998
$data[106,2] = 6.617085102115009
$data[107,0] = 107
$data[107,1] = 117.6999999999998
$data[107,2] = -15.73650902978004
$data[108,0] = 108
$data[108,1] = 118.7999999999998
$data[108,2] = -4.224035879589111
$data[109,0] = 109
$data[109,1] = 119.8999999999998
$data[109,2] = -17.94925916833105
$data[110,0] = 110
$data[110,1] = 120.9999999999998
$data[110,2] = 2.597982777644918
$data[111,0] = 111
$data[111,1] = 122.0999999999998
$data[111,2] = 14.54918396463502

$ws.Range("A2:C113").Value2 = $data

# New rows (86..113) need the same style as the existing index column (A2:A85, style index 1:
# centered/top-aligned, thin border, bold) -- copy formats from A85 down into A86:A113.
$ws.Range("A85").Copy()
$ws.Range("A86:A113").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Updated A2:C113 with new plot data"
